$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 86676.914
$ws.Range("I28").Value = 92828.82000000001
$ws.Range("K28").Value = 92828.82000000001
$ws.Range("M28").Value = -92343.82000000001

$ws.Range("H107").Value = 2485.9285
$ws.Range("I107").Value = 2233.5833
$ws.Range("J107").Value = 4000
$ws.Range("K107").Value = 2233.5833
$ws.Range("L107").Value = 4000
$ws.Range("M107").Value = -313.5832999999998
$ws.Range("N107").Value = -7840

$ws.Range("H113").Value = 9145.083000000001
$ws.Range("I113").Value = 5996.5
$ws.Range("J113").Value = 10719.375
$ws.Range("K113").Value = 5996.5
$ws.Range("L113").Value = 10719.375
$ws.Range("M113").Value = -2742.5
$ws.Range("N113").Value = -17227.375

$ws.Range("H116").Value = 9731.333000000001
$ws.Range("J116").Value = 10900
$ws.Range("L116").Value = 10900
$ws.Range("N116").Value = -17784

$ws.Range("H133").Value = 60649.5
$ws.Range("J133").Value = 60649.5
$ws.Range("L133").Value = 60649.5
$ws.Range("N133").Value = -70769.5

$ws.Range("H137").Value = 3658.25
$ws.Range("I137").Value = 1949.3334
$ws.Range("K137").Value = 5848.0002
$ws.Range("M137").Value = -3298.0002

$ws.Range("H141").Value = 6421.143
$ws.Range("I141").Value = 3494.3333
$ws.Range("K141").Value = 10482.9999
$ws.Range("M141").Value = -5302.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5086.875
$ws.Range("I61").Value = 5086.875
$ws.Range("K61").Value = 5086.875
$ws.Range("M61").Value = -4874.875

$ws.Range("H97").Value = 2950
$ws.Range("I97").Value = 3000
$ws.Range("J97").Value = 2900
$ws.Range("K97").Value = 3000
$ws.Range("L97").Value = 2900
$ws.Range("M97").Value = -2504
$ws.Range("N97").Value = -3892

$ws.Range("H102").Value = 2576.3333
$ws.Range("I102").Value = 2576.3333
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2576.3333
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -954.3332999999998
$ws.Range("N102").ClearContents()

$ws.Range("H122").Value = 2625.5881
$ws.Range("I122").Value = 1847.88
$ws.Range("J122").Value = 4785.8887
$ws.Range("K122").Value = 5543.64
$ws.Range("L122").Value = 14357.6661
$ws.Range("M122").Value = -3093.64
$ws.Range("N122").Value = -19257.6661

$ws.Range("H132").Value = 2837.96
$ws.Range("I132").Value = 807.4103
$ws.Range("K132").Value = 2422.2309
$ws.Range("M132").Value = 107.7691

$ws.Range("H136").Value = 5086.875
$ws.Range("I136").Value = 5086.875
$ws.Range("K136").Value = 15260.625
$ws.Range("M136").Value = -12710.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 3054.6
$ws.Range("I64").Value = 1574.6666
$ws.Range("K64").Value = 1574.6666
$ws.Range("M64").Value = -1349.6666

$ws.Range("H67").Value = 3054.6
$ws.Range("I67").Value = 1574.6666
$ws.Range("K67").Value = 1574.6666
$ws.Range("M67").Value = -794.6666

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1433.6364
$ws.Range("I16").Value = 1339.1428
$ws.Range("J16").Value = 1599
$ws.Range("K16").Value = 1339.1428
$ws.Range("L16").Value = 1599
$ws.Range("M16").Value = -1052.1428
$ws.Range("N16").Value = -2173

$ws.Range("H31").Value = 41040.22
$ws.Range("J31").Value = 64564.79
$ws.Range("L31").Value = 64564.79
$ws.Range("N31").Value = -65154.79

$ws.Range("H34").Value = 41040.22
$ws.Range("J34").Value = 64564.79
$ws.Range("L34").Value = 64564.79
$ws.Range("N34").Value = -64968.79

$ws.Range("H113").Value = 1433.6364
$ws.Range("I113").Value = 1339.1428
$ws.Range("J113").Value = 1599
$ws.Range("K113").Value = 1339.1428
$ws.Range("L113").Value = 1599
$ws.Range("M113").Value = 830.8571999999999
$ws.Range("N113").Value = -5939

$ws.Range("H122").Value = 5826.9443
$ws.Range("J122").Value = 12474.167
$ws.Range("L122").Value = 37422.501
$ws.Range("N122").Value = -42322.501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 78.666664
$ws.Range("J12").Value = 86
$ws.Range("L12").Value = 258
$ws.Range("N12").Value = -604

$ws.Range("H113").Value = 1322
$ws.Range("I113").Value = 1398
$ws.Range("K113").Value = 4194
$ws.Range("M113").Value = -2024

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1804.6666
$ws.Range("I107").Value = 780.25
$ws.Range("K107").Value = 780.25
$ws.Range("M107").Value = 1139.75

$ws.Range("H113").Value = 9997.111000000001
$ws.Range("I113").Value = 9167.5
$ws.Range("K113").Value = 9167.5
$ws.Range("M113").Value = -6997.5

$ws.Range("H132").Value = 41665.242
$ws.Range("I132").Value = 52311.316
$ws.Range("K132").Value = 156933.948
$ws.Range("M132").Value = -154403.948

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 11068.934
$ws.Range("I61").Value = 9534
$ws.Range("J61").Value = 14138.8
$ws.Range("K61").Value = 9534
$ws.Range("L61").Value = 14138.8
$ws.Range("M61").Value = -9332
$ws.Range("N61").Value = -14542.8

$ws.Range("H113").Value = 11068.934
$ws.Range("I113").Value = 9534
$ws.Range("J113").Value = 14138.8
$ws.Range("K113").Value = 9534
$ws.Range("L113").Value = 14138.8
$ws.Range("M113").Value = -7364
$ws.Range("N113").Value = -18478.8

$ws.Range("H132").Value = 3070.32
$ws.Range("I132").Value = 1238.3125
$ws.Range("J132").Value = 6327.222
$ws.Range("K132").Value = 3714.9375
$ws.Range("L132").Value = 18981.666
$ws.Range("M132").Value = -1184.9375
$ws.Range("N132").Value = -24041.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 50000
$ws.Range("I82").Value = 50000
$ws.Range("K82").Value = 50000
$ws.Range("M82").Value = -49617

$ws.Range("H85").Value = 50000
$ws.Range("I85").Value = 50000
$ws.Range("K85").Value = 50000
$ws.Range("M85").Value = -48674

$ws.Range("H122").Value = 8546.379000000001
$ws.Range("I122").Value = 2677.1875
$ws.Range("K122").Value = 8031.5625
$ws.Range("M122").Value = -5581.5625

$ws.Range("H132").Value = 4696.3447
$ws.Range("I132").Value = 5037.4546
$ws.Range("J132").Value = 3624.2856
$ws.Range("K132").Value = 15112.3638
$ws.Range("L132").Value = 10872.8568
$ws.Range("M132").Value = -12582.3638
$ws.Range("N132").Value = -15932.8568

$ws.Range("H136").Value = 5573.519
$ws.Range("I136").Value = 4000.3914
$ws.Range("J136").Value = 17634.166
$ws.Range("K136").Value = 12001.1742
$ws.Range("L136").Value = 52902.49800000001
$ws.Range("M136").Value = -9451.174199999999
$ws.Range("N136").Value = -58002.49800000001
